$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: S.No 29, Date 12-12-2024 (45638)
$ws.Cells.Item(31, 4).Value = 29
$ws.Cells.Item(31, 5).Value = 45638
$ws.Cells.Item(31, 6).Value = "0hr00min"
$ws.Cells.Item(31, 7).Value = "0hr00min"
$ws.Cells.Item(31, 8).Value = "0hr00min"
$ws.Cells.Item(31, 9).Value = "2hr00min"
$ws.Cells.Item(31, 10).Value = "0hr00min"
$ws.Cells.Item(31, 11).Value = "2hr00min"

# Row 32: S.No 30, Date 13-12-2024 (45639)
$ws.Cells.Item(32, 4).Value = 30
$ws.Cells.Item(32, 5).Value = 45639
$ws.Cells.Item(32, 6).Value = "0hr00min"
$ws.Cells.Item(32, 7).Value = "0hr00min"
$ws.Cells.Item(32, 8).Value = "0hr00min"
$ws.Cells.Item(32, 9).Value = "3hr00min"
$ws.Cells.Item(32, 10).Value = "0hr00min"
$ws.Cells.Item(32, 11).Value = "3hr00min"

# Row 33: S.No 31, Date 14-12-2024 (45640), no other data yet
$ws.Cells.Item(33, 4).Value = 31
$ws.Cells.Item(33, 5).Value = 45640

# Copy the date number format (style) from the last existing date cell (E30)
# onto the three newly added date cells so they share the same style index.
$ws.Range("E30").Copy()
$ws.Range("E31:E33").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("F33").Select()
